# Architecture.pptx update — "Processed" stage renamed to "Cleansed"
# (plus the notes-master "last updated" date bump).
#
# NOTE on performance: this deck has many grouped shapes. Looping over a
# COM collection (`.Count` / `.Item`) that is held only inside a function
# parameter is pathologically slow in this host, so every loop below first
# copies the collection/count into local variables before iterating, and
# the whole script is written flat (no recursive shape-search helper).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Notes master "date updated" field: 07/07/16 -> 22/07/16
# ---------------------------------------------------------------------
$nm = $p.NotesMaster
$dt = $nm.HeadersFooters.DateAndTime
$dt.Text = "22/07/16"
Write-Host "date field -> $($nm.Shapes.Item(2).TextFrame.TextRange.Text)"

# ---------------------------------------------------------------------
# 2) Slide 5 — "Processed" -> "Cleansed" (4 boxes) and "Process" -> "Cleanse"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5shapes = $s5.Shapes

$s5a = $s5shapes.Item(21).GroupItems.Item(2)
Write-Host "s5/21: $($s5a.TextFrame.TextRange.Text)"
$s5a.TextFrame.TextRange.Text = "Cleansed"

$s5b = $s5shapes.Item(22).GroupItems.Item(2)
Write-Host "s5/22: $($s5b.TextFrame.TextRange.Text)"
$s5b.TextFrame.TextRange.Text = "Cleansed"

$s5c = $s5shapes.Item(25).GroupItems.Item(2)
Write-Host "s5/25: $($s5c.TextFrame.TextRange.Text)"
$s5c.TextFrame.TextRange.Text = "Cleansed"

$s5d = $s5shapes.Item(26).GroupItems.Item(2)
Write-Host "s5/26: $($s5d.TextFrame.TextRange.Text)"
$s5d.TextFrame.TextRange.Text = "Cleansed"

$s5e = $s5shapes.Item(31).GroupItems.Item(2)
Write-Host "s5/31: $($s5e.TextFrame.TextRange.Text)"
$s5e.TextFrame.TextRange.Text = "Cleanse"

# ---------------------------------------------------------------------
# 3) Slide 6 — "Processed" -> "Cleansed" (4 boxes)
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6shapes = $s6.Shapes

$s6a = $s6shapes.Item(13).GroupItems.Item(2)
Write-Host "s6/13: $($s6a.TextFrame.TextRange.Text)"
$s6a.TextFrame.TextRange.Text = "Cleansed"

$s6b = $s6shapes.Item(14).GroupItems.Item(2)
Write-Host "s6/14: $($s6b.TextFrame.TextRange.Text)"
$s6b.TextFrame.TextRange.Text = "Cleansed"

$s6c = $s6shapes.Item(17).GroupItems.Item(2)
Write-Host "s6/17: $($s6c.TextFrame.TextRange.Text)"
$s6c.TextFrame.TextRange.Text = "Cleansed"

$s6d = $s6shapes.Item(18).GroupItems.Item(2)
Write-Host "s6/18: $($s6d.TextFrame.TextRange.Text)"
$s6d.TextFrame.TextRange.Text = "Cleansed"

# ---------------------------------------------------------------------
# 4) Slide 8 — bullet: "...partitioned, processed" -> "...partitioned, Cleansed"
#    Only the trailing word is replaced so the rest of the run (and its
#    formatting) is left completely untouched; this naturally produces the
#    same two-run split ("...partitioned, " / "Cleansed") that real
#    PowerPoint's edit produced.
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$s8shape = $s8.Shapes.Item(2)
$s8tr = $s8shape.TextFrame.TextRange
$s8para = $s8tr.Paragraphs(3, 1)
Write-Host "s8 para3: $($s8para.Text)"
$s8wordStart = $s8para.Start + 61 - 9
$s8word = $s8tr.Characters($s8wordStart, 9)
Write-Host "s8 word: $($s8word.Text)"
$s8word.Text = "Cleansed"
Write-Host "s8 para3 after: $($s8tr.Paragraphs(3,1).Text)"
